# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) on the per-job leve-profit sheets with
# newly pulled pricing data. Only numeric value cells change; where a
# profit cell no longer applies (no HQ, or price data now unavailable)
# the cell is cleared entirely rather than written as 0, matching source.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 3289.4524
$ws.Range("I43").Value = 1617.125
$ws.Range("K43").Value = 1617.125
$ws.Range("M43").Value = -1548.125

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 112
$ws.Range("H112").Value = 1398.7609
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1398.7609
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4196.2827
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6412.2827

# Row 141
$ws.Range("H141").Value = 915.72
$ws.Range("I141").Value = 915.72
$ws.Range("K141").Value = 2747.16
$ws.Range("M141").Value = 2432.84

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 372.0909
$ws.Range("I4").Value = 293.66666
$ws.Range("J4").Value = 401.5
$ws.Range("K4").Value = 293.66666
$ws.Range("L4").Value = 401.5
$ws.Range("M4").Value = -177.66666
$ws.Range("N4").Value = -633.5

# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 500
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -724

# Row 32
$ws.Range("H32").Value = 15395.647
$ws.Range("I32").Value = 16907
$ws.Range("J32").Value = 11768.4
$ws.Range("K32").Value = 16907
$ws.Range("L32").Value = 11768.4
$ws.Range("M32").Value = -16620
$ws.Range("N32").Value = -12342.4

# Row 61
$ws.Range("H61").Value = 12822409
$ws.Range("I61").Value = 16130842
$ws.Range("J61").Value = 2230.125
$ws.Range("K61").Value = 16130842
$ws.Range("L61").Value = 2230.125
$ws.Range("M61").Value = -16130630
$ws.Range("N61").Value = -2654.125

# Row 74
$ws.Range("H74").Value = 11907877
$ws.Range("I74").Value = 15627153
$ws.Range("J74").Value = 6192.8
$ws.Range("K74").Value = 15627153
$ws.Range("L74").Value = 6192.8
$ws.Range("M74").Value = -15626279
$ws.Range("N74").Value = -7940.8

# Row 77
$ws.Range("H77").Value = 11907877
$ws.Range("I77").Value = 15627153
$ws.Range("J77").Value = 6192.8
$ws.Range("K77").Value = 78135765
$ws.Range("L77").Value = 30964
$ws.Range("M77").Value = -78131397
$ws.Range("N77").Value = -39700

# Row 132
$ws.Range("H132").Value = 5321558
$ws.Range("I132").Value = 11366254
$ws.Range("J132").Value = 2225.56
$ws.Range("K132").Value = 34098762
$ws.Range("L132").Value = 6676.68
$ws.Range("M132").Value = -34096232
$ws.Range("N132").Value = -11736.68

# Row 136
$ws.Range("H136").Value = 12822409
$ws.Range("I136").Value = 16130842
$ws.Range("J136").Value = 2230.125
$ws.Range("K136").Value = 48392526
$ws.Range("L136").Value = 6690.375
$ws.Range("M136").Value = -48389976
$ws.Range("N136").Value = -11790.375

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 500
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -730

# Row 7
$ws.Range("H7").Value = 18597.166
$ws.Range("I7").Value = 25520.75
$ws.Range("J7").Value = 4750
$ws.Range("K7").Value = 25520.75
$ws.Range("L7").Value = 4750
$ws.Range("M7").Value = -25407.75
$ws.Range("N7").Value = -4976

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1461.0303
$ws.Range("I58").Value = 697.7
$ws.Range("J58").Value = 2635.3845
$ws.Range("K58").Value = 697.7
$ws.Range("L58").Value = 2635.3845
$ws.Range("M58").Value = -494.7
$ws.Range("N58").Value = -3041.3845

# Row 99
$ws.Range("H99").Value = 2388.4
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 2980.6667
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 2980.6667
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -5976.6667

# Row 126
$ws.Range("H126").Value = 2388.4
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2980.6667
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 8942.000100000001
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -13882.0001

# Row 132
$ws.Range("H132").Value = 8198371
$ws.Range("I132").Value = 9805329
$ws.Range("J132").Value = 2882.6
$ws.Range("K132").Value = 29415987
$ws.Range("L132").Value = 8647.799999999999
$ws.Range("M132").Value = -29413457
$ws.Range("N132").Value = -13707.8

# Row 134
$ws.Range("H134").Value = 362284.75
$ws.Range("I134").Value = 1311.2128
$ws.Range("K134").Value = 3933.6384
$ws.Range("M134").Value = -1398.6384

# Row 136
$ws.Range("H136").Value = 1461.0303
$ws.Range("I136").Value = 697.7
$ws.Range("J136").Value = 2635.3845
$ws.Range("K136").Value = 2093.1
$ws.Range("L136").Value = 7906.1535
$ws.Range("M136").Value = 456.8999999999996
$ws.Range("N136").Value = -13006.1535

$ws = $wb.Worksheets.Item("CUL")
# Row 118
$ws.Range("H118").Value = 1303.15
$ws.Range("I118").Value = 289
$ws.Range("J118").Value = 1356.5264
$ws.Range("K118").Value = 867
$ws.Range("L118").Value = 4069.5792
$ws.Range("M118").Value = 376
$ws.Range("N118").Value = -6555.5792

# Row 137
$ws.Range("H137").Value = 3973251.2
$ws.Range("I137").Value = 9261540
$ws.Range("J137").Value = 7034.4585
$ws.Range("K137").Value = 27784620
$ws.Range("L137").Value = 21103.3755
$ws.Range("M137").Value = -27779520
$ws.Range("N137").Value = -31303.3755

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 41.727272
$ws.Range("I2").Value = 52
$ws.Range("J2").Value = 29.4
$ws.Range("K2").Value = 52
$ws.Range("L2").Value = 29.4
$ws.Range("M2").Value = 61
$ws.Range("N2").Value = -255.4

# Row 57
$ws.Range("H57").Value = 12224.4
$ws.Range("I57").Value = 3500
$ws.Range("J57").Value = 18040.666
$ws.Range("K57").Value = 3500
$ws.Range("L57").Value = 18040.666
$ws.Range("M57").Value = -2680
$ws.Range("N57").Value = -19680.666

# Row 69
$ws.Range("H69").Value = 35201
$ws.Range("J69").Value = 35201
$ws.Range("L69").Value = 35201
$ws.Range("N69").Value = -36699

# Row 72
$ws.Range("H72").Value = 35201
$ws.Range("J72").Value = 35201
$ws.Range("L72").Value = 105603
$ws.Range("N72").Value = -113091

# Row 122
$ws.Range("H122").Value = 4447064.5
$ws.Range("I122").Value = 7409684.5
$ws.Range("K122").Value = 22229053.5
$ws.Range("M122").Value = -22226603.5

# Row 132
$ws.Range("H132").Value = 3115.4717
$ws.Range("I132").Value = 2159.7144
$ws.Range("K132").Value = 6479.1432
$ws.Range("M132").Value = -3949.1432

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6663.2607
$ws.Range("I7").Value = 12990.8
$ws.Range("J7").Value = 4905.6113
$ws.Range("K7").Value = 12990.8
$ws.Range("L7").Value = 4905.6113
$ws.Range("M7").Value = -12878.8
$ws.Range("N7").Value = -5129.6113

# Row 126
$ws.Range("H126").Value = 6663.2607
$ws.Range("I126").Value = 12990.8
$ws.Range("J126").Value = 4905.6113
$ws.Range("K126").Value = 38972.39999999999
$ws.Range("L126").Value = 14716.8339
$ws.Range("M126").Value = -36502.39999999999
$ws.Range("N126").Value = -19656.8339

# Row 136
$ws.Range("H136").Value = 8623582
$ws.Range("I136").Value = 9616399
$ws.Range("J136").Value = 19167.5
$ws.Range("K136").Value = 28849197
$ws.Range("L136").Value = 57502.5
$ws.Range("M136").Value = -28846647
$ws.Range("N136").Value = -62602.5

$ws = $wb.Worksheets.Item("WVR")
# Row 108
$ws.Range("H108").Value = 31626
$ws.Range("J108").Value = 31626
$ws.Range("L108").Value = 31626
$ws.Range("N108").Value = -39306

# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 126
$ws.Range("H126").Value = 3046.1177
$ws.Range("I126").Value = 2027.0714
$ws.Range("J126").Value = 7801.6665
$ws.Range("K126").Value = 6081.2142
$ws.Range("L126").Value = 23404.9995
$ws.Range("M126").Value = -3611.2142
$ws.Range("N126").Value = -28344.9995

# Row 132
$ws.Range("H132").Value = 1432.877
$ws.Range("I132").Value = 1116.4736
$ws.Range("J132").Value = 3687.25
$ws.Range("K132").Value = 3349.4208
$ws.Range("L132").Value = 11061.75
$ws.Range("M132").Value = -819.4207999999999
$ws.Range("N132").Value = -16121.75

# Row 136
$ws.Range("H136").Value = 766.5538299999999
$ws.Range("I136").Value = 592.46155
$ws.Range("J136").Value = 1462.9231
$ws.Range("K136").Value = 1777.38465
$ws.Range("L136").Value = 4388.7693
$ws.Range("M136").Value = 772.61535
$ws.Range("N136").Value = -9488.7693
